$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.353.86"
$ws.Range("E2").Value = "  +3.45%  "
# Row 3
$ws.Range("D3").Value = "2.609.83"
$ws.Range("E3").Value = "  +1.15%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.34"
$ws.Range("E5").Value = "  -0.18%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.61"
$ws.Range("E6").Value = "  +0.22%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.08%  "
# Row 8
$ws.Range("E8").Value = "  +0.99%  "
# Row 9
$ws.Range("D9").Value = "2.635.16"
$ws.Range("E9").Value = "  +1.85%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.51"
$ws.Range("E10").Value = "  -2.78%  "
# Row 11
$ws.Range("E11").Value = "  +3.10%  "
# Row 12
$ws.Range("E12").Value = "  -3.29%  "
# Row 13
$ws.Range("E13").Value = "  +6.65%  "
# Row 14
$ws.Range("D14").Value = "3.072.96"
$ws.Range("E14").Value = "  +1.26%  "
# Row 15
$ws.Range("D15").Value = "61.268.00"
$ws.Range("E15").Value = "  +3.31%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.52"
$ws.Range("E16").Value = "  +4.63%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000142"
$ws.Range("E17").Value = "  +3.00%  "
# Row 18
$ws.Range("D18").Value = "2.621.67"
$ws.Range("E18").Value = "  +1.41%  "
# Row 19
$ws.Range("E19").Value = "  +3.22%  "
# Row 20
$ws.Range("E20").Value = "  +9.19%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.39"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.19"
$ws.Range("E22").Value = "  +14.43%  "
# Row 23
$ws.Range("E23").Value = "  +0.17%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.518"
$ws.Range("E24").Value = "  +12.64%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.47"
$ws.Range("E25").Value = "  -0.28%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.43%  "
# Row 27
$ws.Range("E27").Value = "  +0.09%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  +5.83%  "
# Row 29
$ws.Range("D29").Value = "0.0₃0798"
$ws.Range("E29").Value = "  +1.99%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("E30").Value = "  +7.56%  "
# Row 31
$ws.Range("E31").Value = "  -0.07%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.36"
$ws.Range("E32").Value = "  +4.78%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.11"
$ws.Range("E33").Value = "  +1.92%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.54"
$ws.Range("E34").Value = "  +2.58%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.29"
$ws.Range("E35").Value = "  +6.09%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.966"
$ws.Range("E36").Value = "  +10.06%  "
# Row 37
$ws.Range("E37").Value = "  +4.06%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  +6.71%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.80"
$ws.Range("E39").Value = "  +1.42%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.854"
$ws.Range("E40").Value = "  -2.15%  "
# Row 41
$ws.Range("E41").Value = "  +3.56%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "298.83"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.17"
$ws.Range("E43").Value = "  +8.07%  "
# Row 44
$ws.Range("E44").Value = "  +1.26%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.09%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.606"
$ws.Range("E46").Value = "  +1.98%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0550"
$ws.Range("E47").Value = "  +2.51%  "
# Row 48
$ws.Range("E48").Value = "  +3.39%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.70"
$ws.Range("E49").Value = "  +0.37%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.68"
$ws.Range("E50").Value = "  +5.85%  "
# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.053.32"
$ws.Range("E51").Value = "  +4.95%  "
